$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "lecture" track anchor date (B2) and the "lab" track anchor
# date (B22) forward by 182 days (fall 2017 -> spring 2018). All the
# other date cells in column B are formulas relative to these two
# anchors, so they recalculate automatically.
$ws.Range("B2").Value = 43185
$ws.Range("B22").Value = 43187

# Row 20 becomes "Memorial day (no class)" (no longer a lecture day),
# and row 21 drops the "(part II)" suffix since there's now only one
# "Building Shiny applications" session.
$ws.Range("D20").Value = "Memorial day (no class)"
$ws.Range("C20").Value = $false
$ws.Range("D21").Value = "Building Shiny applications"

# Update the active selection to match the saved workbook state.
$ws.Range("D21").Select() | Out-Null
